$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs","Dcn","Egfr","ECs",2,0.6666666666666666,6.615074666666668,19.845224,0.0008916467884469992,0.0008916467884469989,3,1,2.08532,6.25596,0.01753772176136817,0.01753772176136816,13.79454750389333,124.15092753504,0.00001563745328520097,0.00001563745328520097),
  @("ECs","Dcn","Egfr","FAPs",2,0.6666666666666666,6.615074666666668,19.845224,0.0008916467884469992,0.0008916467884469989,3,1,101.898173,305.694519,0.8569724579756384,0.8569724579756383,674.0640227919174,6066.576205127257,0.0007641167399415089,0.0007641167399415086),
  @("ECs","Dcn","Egfr","sCs",2,0.6666666666666666,6.615074666666668,19.845224,0.0008916467884469992,0.0008916467884469989,3,1,14.921347,44.76404100000001,0.1254898202629935,0.1254898202629935,98.7058245322427,888.3524207901842,0.0001118925952202893,0.0001118925952202893),
  @("FAPs","Dcn","Egfr","ECs",3,1,7285.701009,21857.103027,0.9820406017477925,0.9820406017477923,3,1,2.08532,6.25596,0.01753772176136817,0.01753772176136816,15193.01802808788,136737.1622527909,0.01722275483181935,0.01722275483181934),
  @("FAPs","Dcn","Egfr","FAPs",3,1,7285.701009,21857.103027,0.9820406017477925,0.9820406017477923,3,1,101.898173,305.694519,0.8569724579756384,0.8569724579756383,742399.6218413566,6681596.59657221,0.8415817483116808,0.8415817483116804),
  @("FAPs","Dcn","Egfr","sCs",3,1,7285.701009,21857.103027,0.9820406017477925,0.9820406017477923,3,1,14.921347,44.76404100000001,0.1254898202629935,0.1254898202629935,108712.4728935391,978412.2560418523,0.1232360986042925,0.1232360986042924),
  @("sCs","Dcn","Egfr","ECs",3,1,126.6246363333333,379.873909,0.01706775146376063,0.01706775146376063,3,1,2.08532,6.25596,0.01753772176136817,0.01753772176136816,264.0528866386267,2376.47597974764,0.0002993294762636181,0.000299329476263618),
  @("sCs","Dcn","Egfr","FAPs",3,1,126.6246363333333,379.873909,0.01706775146376063,0.01706775146376063,3,1,101.898173,305.694519,0.8569724579756384,0.8569724579756383,12902.81909915609,116125.3718924048,0.01462659292401625,0.01462659292401624),
  @("sCs","Dcn","Egfr","sCs",3,1,126.6246363333333,379.873909,0.01706775146376063,0.01706775146376063,3,1,14.921347,44.76404100000001,0.1254898202629935,0.1254898202629935,1889.410137478475,17004.69123730627,0.002141829063480766,0.002141829063480765)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
  $row = $rows[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
  }
}